$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B4 (FilesTab) keeps the "File Name" query text (re-set explicitly to be safe).
$ws.Range("B4").Value2 = "SELECT`r`n    f1.file_name AS `"File Name`",`r`n    s.study_name AS `"Study Name`",`r`n    s.phs_accession AS `"Accession`",`r`n    sp.participant_id AS `"Participant Id`",`r`n    COALESCE((`r`n        SELECT`r`n            REPLACE(GROUP_CONCAT(CASE WHEN rn <= 5 THEN smp.sample_id ELSE NULL END, ', '), ', , ', ', ') ||`r`n            CASE WHEN MAX(rn) > 5 THEN ', ...' ELSE '' END`r`n        FROM (`r`n            SELECT`r`n                smp.sample_id,`r`n                ROW_NUMBER() OVER (ORDER BY smp.sample_id) AS rn`r`n            FROM df_sample smp`r`n            WHERE smp.`"participant.study_participant_id`" = sp.study_participant_id`r`n        ) smp`r`n    ), '') AS `"Sample Id`",`r`n    f1.file_type AS `"File Type`",`r`n    gi.library_strategy AS `"Library Strategy`"`r`nFROM `r`n    df_study s`r`nINNER JOIN `r`n    df_participant sp ON sp.`"study.phs_accession`" = s.phs_accession`r`nINNER JOIN  `r`n    df_sample smp ON smp.`"participant.study_participant_id`" = sp.study_participant_id`r`nINNER JOIN `r`n    df_file f1 ON f1.`"sample.sample_id`" = smp.sample_id`r`nINNER JOIN`r`n    df_genomic_info gi ON gi.`"file.file_id`" = f1.file_id`r`nINNER JOIN`r`n    df_diagnosis d ON d.`"participant.study_participant_id`" = sp.study_participant_id`r`nINNER JOIN`r`n    df_program p ON p.program_acronym = s.`"program.program_acronym`"`r`nWHERE `r`n    s.phs_accession = 'phs001524' AND f1.file_type = 'CRAM'`r`nGROUP BY`r`n    f1.file_name,`r`n    s.study_name,`r`n    s.phs_accession,`r`n    sp.participant_id,`r`n    f1.file_type,`r`n    gi.library_strategy`r`nORDER BY `r`n    f1.file_name ASC`r`nLIMIT 100;"

# B3 (SamplesTab) gets the "Sample ID" query text with the Tumor / Analyte Type
# columns removed (CDS "All studies" testcase).
$ws.Range("B3").Value2 = "SELECT`r`n    DISTINCT (smp.sample_id) AS `"Sample ID`",`r`n    sp.participant_id AS `"Participant ID`", `r`n    s.study_name AS `"Study Name`",`r`n    s.phs_accession AS Accession`r`nFROM `r`n    df_participant sp`r`nJOIN `r`n    df_study s ON sp.`"study.phs_accession`" = s.phs_accession`r`nJOIN `r`n    df_sample smp ON smp.`"participant.study_participant_id`" = sp.study_participant_id`r`nJOIN`r`n    df_diagnosis d ON d.`"participant.study_participant_id`" = sp.study_participant_id`r`nJOIN`r`n    df_program p ON p.program_acronym = s.`"program.program_acronym`"`r`nJOIN`r`n    df_file f1 ON f1.`"sample.sample_id`" = smp.sample_id`r`nJOIN`r`n    df_genomic_info gi ON gi.`"file.file_id`" = f1.file_id`r`nWHERE `r`n   s.phs_accession = 'phs001524' AND f1.file_type = 'CRAM'`r`nORDER BY `r`n    smp.sample_id ASC`r`nLIMIT 100;"

# Update the selection/scroll position to match the new active cell (B3).
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
